$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 834.75
$ws.Range("I28").Value = 651.8
$ws.Range("K28").Value = 651.8
$ws.Range("M28").Value = -166.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 21433284
$ws.Range("I40").Value = 7377.1665
$ws.Range("J40").Value = 37502716
$ws.Range("K40").Value = 7377.1665
$ws.Range("L40").Value = 37502716
$ws.Range("M40").Value = -7202.1665
$ws.Range("N40").Value = -37503066

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 908.6667
$ws.Range("I41").Value = 687.4
$ws.Range("J41").Value = 1185.25
$ws.Range("K41").Value = 687.4
$ws.Range("L41").Value = 1185.25
$ws.Range("M41").Value = -247.4
$ws.Range("N41").Value = -2065.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7976.278
$ws.Range("J51").Value = 8072.8066
$ws.Range("L51").Value = 8072.8066
$ws.Range("N51").Value = -9040.8066

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2727
$ws.Range("I107").Value = 2721.889
$ws.Range("K107").Value = 2721.889
$ws.Range("M107").Value = -801.8890000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1622685.1
$ws.Range("J112").Value = 1885554.6
$ws.Range("L112").Value = 5656663.800000001
$ws.Range("N112").Value = -5658879.800000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2524.5833
$ws.Range("I132").Value = 2181.5151
$ws.Range("K132").Value = 6544.5453
$ws.Range("M132").Value = -4014.5453

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3371863.8
$ws.Range("I137").Value = 5056199
$ws.Range("J137").Value = 3192.8
$ws.Range("K137").Value = 15168597
$ws.Range("L137").Value = 9578.400000000001
$ws.Range("M137").Value = -15166047
$ws.Range("N137").Value = -14678.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3328.4167
$ws.Range("J138").Value = 2918.125
$ws.Range("L138").Value = 8754.375
$ws.Range("N138").Value = -19034.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 36832.668
$ws.Range("I6").Value = 500
$ws.Range("J6").Value = 54999
$ws.Range("K6").Value = 500
$ws.Range("L6").Value = 54999
$ws.Range("M6").Value = -327
$ws.Range("N6").Value = -55345

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1400.5312
$ws.Range("I32").Value = 1415.8914
$ws.Range("K32").Value = 1415.8914
$ws.Range("M32").Value = -1128.8914

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5808.684
$ws.Range("I61").Value = 2426.5
$ws.Range("K61").Value = 2426.5
$ws.Range("M61").Value = -2214.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1894.6586
$ws.Range("I110").Value = 1800.7428
$ws.Range("J110").Value = 2442.5
$ws.Range("K110").Value = 1800.7428
$ws.Range("L110").Value = 2442.5
$ws.Range("M110").Value = 244.2572
$ws.Range("N110").Value = -6532.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5808.684
$ws.Range("I136").Value = 2426.5
$ws.Range("K136").Value = 7279.5
$ws.Range("M136").Value = -4729.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8107.625
$ws.Range("I134").Value = 8107.625
$ws.Range("K134").Value = 24322.875
$ws.Range("M134").Value = -21787.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 21899.166
$ws.Range("I6").Value = 43434.332
$ws.Range("J6").Value = 364
$ws.Range("K6").Value = 43434.332
$ws.Range("L6").Value = 364
$ws.Range("M6").Value = -43321.332
$ws.Range("N6").Value = -590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1232.4642
$ws.Range("I16").Value = 725.8095
$ws.Range("J16").Value = 2752.4285
$ws.Range("K16").Value = 725.8095
$ws.Range("L16").Value = 2752.4285
$ws.Range("M16").Value = -438.8095
$ws.Range("N16").Value = -3326.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 295554.84
$ws.Range("J31").Value = 1081.0588
$ws.Range("L31").Value = 1081.0588
$ws.Range("N31").Value = -1671.0588

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 295554.84
$ws.Range("J34").Value = 1081.0588
$ws.Range("L34").Value = 1081.0588
$ws.Range("N34").Value = -1485.0588

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1837
$ws.Range("I58").Value = 1828
$ws.Range("J58").Value = 1900
$ws.Range("K58").Value = 1828
$ws.Range("L58").Value = 1900
$ws.Range("M58").Value = -1625
$ws.Range("N58").Value = -2306

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 40562.5
$ws.Range("J68").Value = 40562.5
$ws.Range("L68").Value = 40562.5
$ws.Range("N68").Value = -42060.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 40562.5
$ws.Range("J71").Value = 40562.5
$ws.Range("L71").Value = 121687.5
$ws.Range("N71").Value = -129175.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4163.7407
$ws.Range("I105").Value = 1101.238
$ws.Range("K105").Value = 1101.238
$ws.Range("M105").Value = 645.7619999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 3548.7737
$ws.Range("I107").Value = 919
$ws.Range("J107").Value = 4586.8423
$ws.Range("K107").Value = 919
$ws.Range("L107").Value = 4586.8423
$ws.Range("M107").Value = 1001
$ws.Range("N107").Value = -8426.8423

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1232.4642
$ws.Range("I113").Value = 725.8095
$ws.Range("J113").Value = 2752.4285
$ws.Range("K113").Value = 725.8095
$ws.Range("L113").Value = 2752.4285
$ws.Range("M113").Value = 1444.1905
$ws.Range("N113").Value = -7092.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2583.64
$ws.Range("I134").Value = 2318
$ws.Range("K134").Value = 6954
$ws.Range("M134").Value = -4419

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1837
$ws.Range("I136").Value = 1828
$ws.Range("J136").Value = 1900
$ws.Range("K136").Value = 5484
$ws.Range("L136").Value = 5700
$ws.Range("M136").Value = -2934
$ws.Range("N136").Value = -10800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 572.5
$ws.Range("I34").Value = 101.92857
$ws.Range("J34").Value = 2219.5
$ws.Range("K34").Value = 305.78571
$ws.Range("L34").Value = 6658.5
$ws.Range("M34").Value = -221.78571
$ws.Range("N34").Value = -6826.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2689855.5
$ws.Range("I68").Value = 16667325
$ws.Range("K68").Value = 50001975
$ws.Range("M68").Value = -50001164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 12
$ws.Range("I70").Value = 12
$ws.Range("K70").Value = 36
$ws.Range("M70").Value = 279

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2689855.5
$ws.Range("I71").Value = 16667325
$ws.Range("K71").Value = 150005925
$ws.Range("M71").Value = -150001869

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 12
$ws.Range("I73").Value = 12
$ws.Range("K73").Value = 36
$ws.Range("M73").Value = 1056

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 817.4286
$ws.Range("J86").Value = 869
$ws.Range("L86").Value = 2607
$ws.Range("N86").Value = -4979

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 817.4286
$ws.Range("J89").Value = 869
$ws.Range("L89").Value = 7821
$ws.Range("N89").Value = -19677

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 864.92
$ws.Range("I113").Value = 736.1429000000001
$ws.Range("J113").Value = 915
$ws.Range("K113").Value = 2208.4287
$ws.Range("L113").Value = 2745
$ws.Range("M113").Value = -38.42870000000039
$ws.Range("N113").Value = -7085

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 1905
$ws.Range("J127").Value = 1757.5
$ws.Range("L127").Value = 5272.5
$ws.Range("N127").Value = -15192.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3740
$ws.Range("J126").Value = 4115
$ws.Range("L126").Value = 12345
$ws.Range("N126").Value = -17285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1820.4
$ws.Range("I2").Value = 1034
$ws.Range("K2").Value = 1034
$ws.Range("M2").Value = -922

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5981.778
$ws.Range("I7").Value = 6113.154
$ws.Range("K7").Value = 6113.154
$ws.Range("M7").Value = -6001.154

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 30001
$ws.Range("J97").Value = 30001
$ws.Range("L97").Value = 30001
$ws.Range("N97").Value = -31983

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5981.778
$ws.Range("I126").Value = 6113.154
$ws.Range("K126").Value = 18339.462
$ws.Range("M126").Value = -15869.462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4081.3635
$ws.Range("I132").Value = 4018.7144
$ws.Range("J132").Value = 4191
$ws.Range("K132").Value = 12056.1432
$ws.Range("L132").Value = 12573
$ws.Range("M132").Value = -9526.143199999999
$ws.Range("N132").Value = -17633

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4005.389
$ws.Range("I136").Value = 4228.4287
$ws.Range("J136").Value = 3224.75
$ws.Range("K136").Value = 12685.2861
$ws.Range("L136").Value = 9674.25
$ws.Range("M136").Value = -10135.2861
$ws.Range("N136").Value = -14774.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
